# Ads1299_defRegs.xlsx edit: disable SRB1, change CH1-3SET default from
# 0x60 to 0xE1 (SRB2 closed instead of open) and CH4-7SET from 0xF1 to
# 0xE1 (SRB1 routing bit cleared), and clear LOFF_SENSN bit2 (0x20 -> 0x00).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Move the active selection from M3 to D3 ---------------------------
$ws.Range("D3").Select()

# --- CH1SET/CH2SET/CH3SET (rows 7-9): D1,D0 become 1,1 (SRB2 closed)   --
# and the whole D:K block switches from the mixed "default" styling to
# the same "Explanatory Text" (centered, italic grey) styling used by
# CH4SET..CH7SET below it.
foreach ($row in 7, 8, 9) {
    $rng = $ws.Range("D" + $row + ":K" + $row)
    $rng.Style = "Explanatory Text"
    $rng.HorizontalAlignment = -4108  # xlCenter

    $ws.Cells.Item($row, 4).Value = 1   # D: PDn1       -> 1
    $ws.Cells.Item($row, 5).Value = 1   # E: GAIN2      -> 1
    $ws.Cells.Item($row, 6).Value = 1   # F: GAIN1      -> 1
    $ws.Cells.Item($row, 7).Value = 0   # G: GAIN0      -> 0
    $ws.Cells.Item($row, 8).Value = 0   # H: SRB2       -> 0
    $ws.Cells.Item($row, 9).Value = 0   # I: MUX2       -> 0
    $ws.Cells.Item($row, 10).Value = 0  # J: MUX1       -> 0
    $ws.Cells.Item($row, 11).Value = 1  # K: MUX0       -> 1
}

# --- CH4SET..CH7SET (rows 10-13): disable SRB1 (column G: 1 -> 0) ------
foreach ($row in 10, 11, 12, 13) {
    $ws.Cells.Item($row, 7).Value = 0   # G: SRB1 -> 0
}

# --- LOFF_SENSN (row 22): clear bit (column F: 1 -> 0) ------------------
$ws.Cells.Item(22, 6).Value = 0
